# Update market-price-derived columns (H-N) on several sheets to refreshed values.
$wb = $excel.ActiveWorkbook

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1922.6364
$ws.Range("I40").Value = 1938.7778
$ws.Range("J40").Value = 1850
$ws.Range("K40").Value = 1938.7778
$ws.Range("L40").Value = 1850
$ws.Range("M40").Value = -1763.7778
$ws.Range("N40").Value = -2200

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1717.5
$ws.Range("J112").Value = 1737.2727
$ws.Range("L112").Value = 5211.8181
$ws.Range("N112").Value = -7427.8181

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 7992.174
$ws.Range("I113").Value = 2666.5454
$ws.Range("J113").Value = 12874
$ws.Range("K113").Value = 2666.5454
$ws.Range("L113").Value = 12874
$ws.Range("M113").Value = 587.4546
$ws.Range("N113").Value = -19382

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1682948.1
$ws.Range("I125").Value = 2433
$ws.Range("J125").Value = 3923635
$ws.Range("K125").Value = 21897
$ws.Range("L125").Value = 35312715
$ws.Range("M125").Value = -19437
$ws.Range("N125").Value = -35317635

# ARM row 22
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 14500.25
$ws.Range("I22").Value = 2667
$ws.Range("J22").Value = 50000
$ws.Range("K22").Value = 2667
$ws.Range("L22").Value = 50000
$ws.Range("M22").Value = -2368
$ws.Range("N22").Value = -50598

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1100.5
$ws.Range("I74").Value = 1007.64703
$ws.Range("J74").Value = 1244
$ws.Range("K74").Value = 1007.64703
$ws.Range("L74").Value = 1244
$ws.Range("M74").Value = -133.64703
$ws.Range("N74").Value = -2992

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1100.5
$ws.Range("I77").Value = 1007.64703
$ws.Range("J77").Value = 1244
$ws.Range("K77").Value = 5038.23515
$ws.Range("L77").Value = 6220
$ws.Range("M77").Value = -670.2351499999995
$ws.Range("N77").Value = -14956

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2551.5386
$ws.Range("I122").Value = 980
$ws.Range("J122").Value = 3898.5715
$ws.Range("K122").Value = 2940
$ws.Range("L122").Value = 11695.7145
$ws.Range("M122").Value = -490
$ws.Range("N122").Value = -16595.7145

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2156.5476
$ws.Range("I132").Value = 1082.3462
$ws.Range("J132").Value = 3902.125
$ws.Range("K132").Value = 3247.0386
$ws.Range("L132").Value = 11706.375
$ws.Range("M132").Value = -717.0385999999999
$ws.Range("N132").Value = -16766.375

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1631.4375
$ws.Range("I86").Value = 1868.9474
$ws.Range("K86").Value = 1868.9474
$ws.Range("M86").Value = -745.9474

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1631.4375
$ws.Range("I89").Value = 1868.9474
$ws.Range("K89").Value = 9344.737000000001
$ws.Range("M89").Value = -3728.737000000001

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2542.7058
$ws.Range("I31").Value = 2465.889
$ws.Range("J31").Value = 2629.125
$ws.Range("K31").Value = 2465.889
$ws.Range("L31").Value = 2629.125
$ws.Range("M31").Value = -2170.889
$ws.Range("N31").Value = -3219.125

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2542.7058
$ws.Range("I34").Value = 2465.889
$ws.Range("J34").Value = 2629.125
$ws.Range("K34").Value = 2465.889
$ws.Range("L34").Value = 2629.125
$ws.Range("M34").Value = -2263.889
$ws.Range("N34").Value = -3033.125

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4160.154
$ws.Range("I132").Value = 3123.75
$ws.Range("J132").Value = 5818.4
$ws.Range("K132").Value = 9371.25
$ws.Range("L132").Value = 17455.2
$ws.Range("M132").Value = -6841.25
$ws.Range("N132").Value = -22515.2

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 200
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 600
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -431
$ws.Range("N17").ClearContents()

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1608.3125
$ws.Range("I102").Value = 1269.4166
$ws.Range("J102").Value = 2625
$ws.Range("K102").Value = 1269.4166
$ws.Range("L102").Value = 2625
$ws.Range("M102").Value = 352.5834
$ws.Range("N102").Value = -5869

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2338.7222
$ws.Range("I122").Value = 2031.125
$ws.Range("J122").Value = 2953.9167
$ws.Range("K122").Value = 6093.375
$ws.Range("L122").Value = 8861.750100000001
$ws.Range("M122").Value = -3643.375
$ws.Range("N122").Value = -13761.7501

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1845.75
$ws.Range("I126").Value = 1577.258
$ws.Range("J126").Value = 2335.353
$ws.Range("K126").Value = 4731.774
$ws.Range("L126").Value = 7006.059
$ws.Range("M126").Value = -2261.774
$ws.Range("N126").Value = -11946.059

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4199.3335
$ws.Range("I132").Value = 1795.75
$ws.Range("J132").Value = 9006.5
$ws.Range("K132").Value = 5387.25
$ws.Range("L132").Value = 27019.5
$ws.Range("M132").Value = -2857.25
$ws.Range("N132").Value = -32079.5

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 51937.65
$ws.Range("I7").Value = 78330.69500000001
$ws.Range("J7").Value = 2922
$ws.Range("K7").Value = 78330.69500000001
$ws.Range("L7").Value = 2922
$ws.Range("M7").Value = -78218.69500000001
$ws.Range("N7").Value = -3146

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 30623.889
$ws.Range("I40").Value = 43600.832
$ws.Range("J40").Value = 4670
$ws.Range("K40").Value = 43600.832
$ws.Range("L40").Value = 4670
$ws.Range("M40").Value = -43464.832
$ws.Range("N40").Value = -4942

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7411858
$ws.Range("I122").Value = 55556656
$ws.Range("J122").Value = 4966.077
$ws.Range("K122").Value = 166669968
$ws.Range("L122").Value = 14898.231
$ws.Range("M122").Value = -166667518
$ws.Range("N122").Value = -19798.231

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 51937.65
$ws.Range("I126").Value = 78330.69500000001
$ws.Range("J126").Value = 2922
$ws.Range("K126").Value = 234992.085
$ws.Range("L126").Value = 8766
$ws.Range("M126").Value = -232522.085
$ws.Range("N126").Value = -13706

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 13254.818
$ws.Range("I132").Value = 24289.223
$ws.Range("J132").Value = 5615.615
$ws.Range("K132").Value = 72867.66900000001
$ws.Range("L132").Value = 16846.845
$ws.Range("M132").Value = -70337.66900000001
$ws.Range("N132").Value = -21906.845

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 15877064
$ws.Range("I136").Value = 4249.5835
$ws.Range("K136").Value = 12748.7505
$ws.Range("M136").Value = -10198.7505

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1007.96
$ws.Range("I113").Value = 1044.4667
$ws.Range("J113").Value = 953.2
$ws.Range("K113").Value = 3133.4001
$ws.Range("L113").Value = 2859.6
$ws.Range("M113").Value = -963.4000999999998
$ws.Range("N113").Value = -7199.6

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 52073.1
$ws.Range("I122").Value = 60910.94
$ws.Range("J122").Value = 1992
$ws.Range("K122").Value = 182732.82
$ws.Range("L122").Value = 5976
$ws.Range("M122").Value = -180282.82
$ws.Range("N122").Value = -10876

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 41285.76
$ws.Range("I126").Value = 63372.375
$ws.Range("J126").Value = 2020.6666
$ws.Range("K126").Value = 190117.125
$ws.Range("L126").Value = 6061.9998
$ws.Range("M126").Value = -187647.125
$ws.Range("N126").Value = -11001.9998

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2031.5312
$ws.Range("I132").Value = 1584.24
$ws.Range("J132").Value = 3629
$ws.Range("K132").Value = 4752.72
$ws.Range("L132").Value = 10887
$ws.Range("M132").Value = -2222.72
$ws.Range("N132").Value = -15947
